# Scheduled runner update: refresh currentAveragePrice / Leve price & profit
# figures across the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6588.3335
$ws.Range("I32").Value = 7960
$ws.Range("J32").Value = 4873.75
$ws.Range("K32").Value = 7960
$ws.Range("L32").Value = 4873.75
$ws.Range("M32").Value = -7634
$ws.Range("N32").Value = -5525.75

$ws.Range("H70").Value = 1596.5
$ws.Range("I70").Value = 1695.5
$ws.Range("K70").Value = 5086.5
$ws.Range("M70").Value = -4816.5

$ws.Range("H73").Value = 1596.5
$ws.Range("I73").Value = 1695.5
$ws.Range("K73").Value = 5086.5
$ws.Range("M73").Value = -4150.5

$ws.Range("H112").Value = 1098.5
$ws.Range("J112").Value = 1098.5
$ws.Range("L112").Value = 3295.5
$ws.Range("N112").Value = -5511.5

$ws.Range("H137").Value = 3799.6
$ws.Range("I137").Value = 1998
$ws.Range("J137").Value = 4250
$ws.Range("K137").Value = 5994
$ws.Range("L137").Value = 12750
$ws.Range("M137").Value = -3444
$ws.Range("N137").Value = -17850

$ws.Range("H138").Value = 7252.9487
$ws.Range("J138").Value = 7878.706
$ws.Range("L138").Value = 23636.118
$ws.Range("N138").Value = -33916.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2628.6667
$ws.Range("I32").Value = 2585.8125
$ws.Range("K32").Value = 2585.8125
$ws.Range("M32").Value = -2298.8125

$ws.Range("H74").Value = 6000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -7748

$ws.Range("H77").Value = 6000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -38736

$ws.Range("H102").Value = 1300.4166
$ws.Range("I102").Value = 1365.6
$ws.Range("K102").Value = 1365.6
$ws.Range("M102").Value = 256.4000000000001

$ws.Range("H110").Value = 4115.5713
$ws.Range("I110").Value = 4222.6
$ws.Range("J110").Value = 3848
$ws.Range("K110").Value = 4222.6
$ws.Range("L110").Value = 3848
$ws.Range("M110").Value = -2177.6
$ws.Range("N110").Value = -7938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5889.3335
$ws.Range("I31").Value = 1829.5
$ws.Range("K31").Value = 1829.5
$ws.Range("M31").Value = -1534.5

$ws.Range("H34").Value = 5889.3335
$ws.Range("I34").Value = 1829.5
$ws.Range("K34").Value = 1829.5
$ws.Range("M34").Value = -1627.5

$ws.Range("H105").Value = 1421.1428
$ws.Range("I105").Value = 1241.5
$ws.Range("K105").Value = 1241.5
$ws.Range("M105").Value = 505.5

$ws.Range("H133").Value = 124800
$ws.Range("J133").Value = 124800
$ws.Range("L133").Value = 124800
$ws.Range("N133").Value = -129860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 90.5
$ws.Range("I7").Value = 81
$ws.Range("K7").Value = 243
$ws.Range("M7").Value = -131

$ws.Range("H40").Value = 200
$ws.Range("I40").Value = 200
$ws.Range("K40").Value = 800
$ws.Range("M40").Value = -731

$ws.Range("H68").Value = 1288.8
$ws.Range("J68").Value = 1311
$ws.Range("L68").Value = 3933
$ws.Range("N68").Value = -5555

$ws.Range("H71").Value = 1288.8
$ws.Range("J71").Value = 1311
$ws.Range("L71").Value = 11799
$ws.Range("N71").Value = -19911

$ws.Range("H80").Value = 9002
$ws.Range("I80").Value = 9002
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 27006
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -26070
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 9002
$ws.Range("I83").Value = 9002
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 81018
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -76338
$ws.Range("N83").Value = ""

$ws.Range("H132").Value = 3634.1177
$ws.Range("J132").Value = 5121.222
$ws.Range("L132").Value = 46090.998
$ws.Range("N132").Value = -51150.998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I126").Value = 2248.8333
$ws.Range("J126").Value = 1998.25
$ws.Range("K126").Value = 6746.499899999999
$ws.Range("L126").Value = 5994.75
$ws.Range("M126").Value = -4276.499899999999
$ws.Range("N126").Value = -10934.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2934.7856
$ws.Range("I40").Value = 2852.4614
$ws.Range("J40").Value = 4005
$ws.Range("K40").Value = 2852.4614
$ws.Range("L40").Value = 4005
$ws.Range("M40").Value = -2716.4614
$ws.Range("N40").Value = -4277

$ws.Range("H46").Value = 2850.6667
$ws.Range("I46").Value = 2850.6667
$ws.Range("K46").Value = 2850.6667
$ws.Range("M46").Value = -2662.6667

$ws.Range("H122").Value = 1664.3334
$ws.Range("I122").Value = 1664.3334
$ws.Range("K122").Value = 4993.0002
$ws.Range("M122").Value = -2543.0002

$ws.Range("H132").Value = 3111.75
$ws.Range("I132").Value = 2149.5833
$ws.Range("J132").Value = 5998.25
$ws.Range("K132").Value = 6448.749899999999
$ws.Range("L132").Value = 17994.75
$ws.Range("M132").Value = -3918.749899999999
$ws.Range("N132").Value = -23054.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13501
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""

$ws.Range("H84").Value = 13501
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""

$ws.Range("H96").Value = 2905.2
$ws.Range("I96").Value = 3137
$ws.Range("K96").Value = 3137
$ws.Range("M96").Value = -1764

$ws.Range("H107").Value = 404.7143
$ws.Range("I107").Value = 404.7143
$ws.Range("K107").Value = 1214.1429
$ws.Range("M107").Value = 705.8571000000002

$ws.Range("H122").Value = 1457.7646
$ws.Range("I122").Value = 1298.9375
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 3896.8125
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -1446.8125
$ws.Range("N122").Value = -16897

$ws.Range("H132").Value = 2651
$ws.Range("I132").Value = 2397.5
$ws.Range("J132").Value = 3348.125
$ws.Range("K132").Value = 7192.5
$ws.Range("L132").Value = 10044.375
$ws.Range("M132").Value = -4662.5
$ws.Range("N132").Value = -15104.375

$ws.Range("H136").Value = 5422.8213
$ws.Range("I136").Value = 6430.8184
$ws.Range("K136").Value = 19292.4552
$ws.Range("M136").Value = -16742.4552
